$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.334.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.32"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.35%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7055"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.23%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "238.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.17%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07746"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.65%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3056"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.44%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +7.41%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08185"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.24%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.890.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.241"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.29%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.7175"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.13%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.09%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.340.31"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.17%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.808"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.53%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "240.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.35%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007806"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.63%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.13%  "

# Row 21
$ws.Range("B21").Value = "Dai"
$ws.Range("C21").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

# Row 22
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.114.06"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.05%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.28%  "

# Row 24
$ws.Range("E24").Value = "  +2.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.06"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.04%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.920"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.74%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1447"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.45%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.14"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.08%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.920"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.65%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.361"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.04%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.516"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.29%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.308"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.85%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.042"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.50%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05212"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.187"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.30%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7146"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.95%  "

# Row 37
$ws.Range("E37").Value = "  -0.15%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.682"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.47%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01848"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.20%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.703"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.86%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.180.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.38%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9140"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.47%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.995"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.34%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.27"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.57%  "

# Row 45
$ws.Range("B45").Value = "TheSandbox"
$ws.Range("C45").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4273"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.29%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.28%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.90"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.74%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5365"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.90%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.753"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.17%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.170"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.13%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.010"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "
